$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the two new header columns BC1/BD1, copying format from BB1
$ws.Range("BB1").Copy()
$ws.Range("BC1:BD1").PasteSpecial(-4122)
$ws.Range("BC1").Value2 = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value2 = "Odd_CS_4-4_HT"

# 2) Insert a new row at position 2, shifting the existing match down to row 3
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Range("A2:BD2").ClearFormats()

# 3) Populate the new row 2 with the new match data
$ws.Range("A2").Value2 = "Ukw2Lkbe"
$ws.Range("B2").Value2 = "24/11/2024"
$ws.Range("C2").Value2 = "06:30"
$ws.Range("D2").Value2 = "AUSTRIA - 2. LIGA"
$ws.Range("E2").Value2 = "A. Lustenau"
$ws.Range("F2").Value2 = "Stripfing"
$ws.Range("G2").Value2 = 1.62
$ws.Range("H2").Value2 = 3.8
$ws.Range("I2").Value2 = 4.9
$ws.Range("J2").Value2 = 2.18
$ws.Range("K2").Value2 = 2.2
$ws.Range("L2").Value2 = 5
$ws.Range("M2").Value2 = 1.05
$ws.Range("N2").Value2 = 7.8
$ws.Range("O2").Value2 = 1.26
$ws.Range("P2").Value2 = 3.5
$ws.Range("Q2").Value2 = 1.78
$ws.Range("R2").Value2 = 1.95
$ws.Range("S2").Value2 = 1.38
$ws.Range("T2").Value2 = 2.82
$ws.Range("U2").Value2 = 1.8
$ws.Range("V2").Value2 = 1.91
$ws.Range("W2").Value2 = 7.2
$ws.Range("X2").Value2 = 7.8
$ws.Range("Y2").Value2 = 8
$ws.Range("Z2").Value2 = 12.5
$ws.Range("AA2").Value2 = 12.5
$ws.Range("AB2").Value2 = 25
$ws.Range("AC2").Value2 = 7.8
$ws.Range("AD2").Value2 = 7.4
$ws.Range("AE2").Value2 = 16
$ws.Range("AF2").Value2 = 70
$ws.Range("AG2").Value2 = 500
$ws.Range("AH2").Value2 = 13.5
$ws.Range("AI2").Value2 = 28
$ws.Range("AJ2").Value2 = 15.5
$ws.Range("AK2").Value2 = 80
$ws.Range("AL2").Value2 = 45
$ws.Range("AM2").Value2 = 50
$ws.Range("AN2").Value2 = 3.5
$ws.Range("AO2").Value2 = 7.9
$ws.Range("AP2").Value2 = 17.5
$ws.Range("AQ2").Value2 = 26
$ws.Range("AR2").Value2 = 55
$ws.Range("AS2").Value2 = 250
$ws.Range("AT2").Value2 = 2.82
$ws.Range("AU2").Value2 = 7.6
$ws.Range("AV2").Value2 = 70
$ws.Range("AW2").Value2 = 6.5
$ws.Range("AX2").Value2 = 28
$ws.Range("AY2").Value2 = 32
$ws.Range("AZ2").Value2 = 175
$ws.Range("BA2").Value2 = 200
$ws.Range("BB2").Value2 = 450
$ws.Range("BC2").Value2 = 81
$ws.Range("BD2").Value2 = 81

$excel.CutCopyMode = $false
